$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column J: "this col will be hidden" spelled across J1:J5
$ws.Range("J1").Value = "this"
$ws.Range("J2").Value = "col"
$ws.Range("J3").Value = "will"
$ws.Range("J4").Value = "be"
$ws.Range("J5").Value = "hidden"

# Column K: "this col remains" spelled across K1:K3
$ws.Range("K1").Value = "this"
$ws.Range("K2").Value = "col"
$ws.Range("K3").Value = "remains"

# Select K4 to match final state
[void]$ws.Range("K4").Select()
